$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a single "last updated" date value that is
# repeated for every data row (rows 2 through 530). Bump it from 45177
# (2023-09-08) to 45178 (2023-09-09) for every row.
$startRow = 2
$endRow = 530

for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
